# mating_event.xlsx — rename the "form_id" setting to "table_id" and add a
# new "properties" sheet (minimal colOrder property) so the app-designer
# pipeline can generate definitions.csv / properties.csv for this form.

$wb = $excel.ActiveWorkbook

# --- settings sheet: rename the "form_id" row to "table_id" -----------------
$settingsSheet = $wb.Worksheets.Item("settings")
$settingsSheet.Range("A2").Value = "table_id"
$settingsSheet.Range("A3").Select() | Out-Null

# --- add the new "properties" sheet, placed after "settings" ---------------
$propsSheet = $wb.Worksheets.Add([Type]::Missing, $settingsSheet)
$propsSheet.Name = "properties"

# header row
$propsSheet.Range("A1").Value = "partition"
$propsSheet.Range("B1").Value = "aspect"
$propsSheet.Range("C1").Value = "key"
$propsSheet.Range("D1").Value = "type"
$propsSheet.Range("E1").Value = "value"

# single data row: Table / default / colOrder / array / [...]
$propsSheet.Range("A2").Value = "Table"
$propsSheet.Range("B2").Value = "default"
$propsSheet.Range("C2").Value = "colOrder"
$propsSheet.Range("D2").Value = "array"
$propsSheet.Range("E2").Value = '["M_FOL_date","M_FOL_B_focal_AnimID","M_time","M_male","M_female","M_fem_swelling","M_fail_flag","M_interference","M_comments","M_consort_flag","M_mate_guard_flag"]'

# the properties sheet becomes the active / selected sheet & cell
$propsSheet.Range("E5").Select() | Out-Null
